# Applies the scheduled-runner market/profit data refresh to the
# "Maduin_Profits" workbook (per-sheet Leve profit tables).
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for the specific Leve rows touched by this data refresh.

$wb = $excel.ActiveWorkbook

# ============ Sheet: ALC ============
$ws = $wb.Worksheets.Item("ALC")

# Row 55 (Leve Item ID 5517)
$ws.Range("H55").Value = 495.18182
$ws.Range("I55").Value = 74.666664
$ws.Range("J55").Value = 999.8
$ws.Range("K55").Value = 74.666664
$ws.Range("L55").Value = 999.8
$ws.Range("M55").Value = 139.333336
$ws.Range("N55").Value = -1427.8

# Row 103 (Leve Item ID 19909)
$ws.Range("H103").Value = 4666
$ws.Range("I103").Value = 3000
$ws.Range("J103").Value = 5776.6665
$ws.Range("K103").Value = 9000
$ws.Range("L103").Value = 17329.9995
$ws.Range("M103").Value = -8414
$ws.Range("N103").Value = -18501.9995

# Row 105 (Leve Item ID 18668)
$ws.Range("H105").Value = 10000
$ws.Range("J105").Value = 10000
$ws.Range("L105").Value = 10000
$ws.Range("N105").Value = -16988

# Row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 1443.7142
$ws.Range("I106").Value = 1443.7142
$ws.Range("K106").Value = 1443.7142
$ws.Range("M106").Value = -812.7141999999999

# Row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 1195
$ws.Range("I111").Value = 961.1111
$ws.Range("K111").Value = 2883.3333
$ws.Range("M111").Value = 183.6667000000002

# Row 115 (Leve Item ID 27957)
$ws.Range("H115").Value = 649
$ws.Range("I115").Value = 649
$ws.Range("K115").Value = 1947
$ws.Range("M115").Value = -380

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 6249.9443
$ws.Range("J116").Value = 6264.647
$ws.Range("L116").Value = 6264.647
$ws.Range("N116").Value = -13148.647

# Row 131 (Leve Item ID 36108)
$ws.Range("H131").Value = 846.5
$ws.Range("I131").Value = 615.8
$ws.Range("K131").Value = 1847.4
$ws.Range("M131").Value = 3192.6

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 5282.6665
$ws.Range("I132").Value = 1468
$ws.Range("K132").Value = 4404
$ws.Range("M132").Value = -1874

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 2461.963
$ws.Range("I137").Value = 1721
$ws.Range("J137").Value = 3150
$ws.Range("K137").Value = 5163
$ws.Range("L137").Value = 9450
$ws.Range("M137").Value = -2613
$ws.Range("N137").Value = -14550

# Row 140 (Leve Item ID 42459)
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 10799
$ws.Range("I141").Value = 12624.5
$ws.Range("J141").Value = 3497
$ws.Range("K141").Value = 37873.5
$ws.Range("L141").Value = 10491
$ws.Range("M141").Value = -32693.5
$ws.Range("N141").Value = -20851


# ============ Sheet: ARM ============
$ws = $wb.Worksheets.Item("ARM")

# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 4995.9653
$ws.Range("I32").Value = 4380.7856
$ws.Range("K32").Value = 4380.7856
$ws.Range("M32").Value = -4093.7856

# Row 92 (Leve Item ID 18050)
$ws.Range("H92").Value = 42516.668
$ws.Range("J92").Value = 42516.668
$ws.Range("L92").Value = 42516.668
$ws.Range("N92").Value = -47508.668

# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 2961.75
$ws.Range("I102").Value = 1339
$ws.Range("K102").Value = 1339
$ws.Range("M102").Value = 283


# ============ Sheet: BSM ============
$ws = $wb.Worksheets.Item("BSM")

# Row 54 (Leve Item ID 2376)
$ws.Range("H54").Value = 5733.1665
$ws.Range("I54").Value = 5733.1665
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 5733.1665
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -5249.1665
$ws.Range("N54").ClearContents()

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 5467.8887
$ws.Range("I105").Value = 4678
$ws.Range("J105").Value = 6099.8
$ws.Range("K105").Value = 4678
$ws.Range("L105").Value = 6099.8
$ws.Range("M105").Value = -2931
$ws.Range("N105").Value = -9593.799999999999


# ============ Sheet: CRP ============
$ws = $wb.Worksheets.Item("CRP")

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 4823.2666
$ws.Range("I31").Value = 3383.7
$ws.Range("K31").Value = 3383.7
$ws.Range("M31").Value = -3088.7

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 4823.2666
$ws.Range("I34").Value = 3383.7
$ws.Range("K34").Value = 3383.7
$ws.Range("M34").Value = -3181.7

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 4966.9443
$ws.Range("I58").Value = 4114.75
$ws.Range("K58").Value = 4114.75
$ws.Range("M58").Value = -3911.75

# Row 95 (Leve Item ID 18192)
$ws.Range("H95").Value = 10812.5
$ws.Range("J95").Value = 10812.5
$ws.Range("L95").Value = 10812.5
$ws.Range("N95").Value = -16304.5

# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 1924.7142
$ws.Range("I105").Value = 1493.25
$ws.Range("K105").Value = 1493.25
$ws.Range("M105").Value = 253.75

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 6832.838
$ws.Range("I132").Value = 5160.567
$ws.Range("K132").Value = 15481.701
$ws.Range("M132").Value = -12951.701

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 4966.9443
$ws.Range("I136").Value = 4114.75
$ws.Range("K136").Value = 12344.25
$ws.Range("M136").Value = -9794.25


# ============ Sheet: CUL ============
$ws = $wb.Worksheets.Item("CUL")

# Row 69 (Leve Item ID 12850)
$ws.Range("H69").Value = 966.3333
$ws.Range("I69").Value = 899
$ws.Range("J69").Value = 1000
$ws.Range("K69").Value = 2697
$ws.Range("L69").Value = 3000
$ws.Range("M69").Value = -1886
$ws.Range("N69").Value = -4622

# Row 72 (Leve Item ID 12850)
$ws.Range("H72").Value = 966.3333
$ws.Range("I72").Value = 899
$ws.Range("K72").Value = 8091
$ws.Range("L72").Value = 9000
$ws.Range("M72").Value = -4035
$ws.Range("N72").Value = -17112

# Row 92 (Leve Item ID 19841)
$ws.Range("H92").Value = 870.125
$ws.Range("I92").Value = 650
$ws.Range("J92").Value = 943.5
$ws.Range("K92").Value = 1950
$ws.Range("L92").Value = 2830.5
$ws.Range("M92").Value = -702
$ws.Range("N92").Value = -5326.5

# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 970.5
$ws.Range("I107").Value = 520.3333
$ws.Range("J107").Value = 1163.4286
$ws.Range("K107").Value = 1560.9999
$ws.Range("L107").Value = 3490.2858
$ws.Range("M107").Value = 359.0001
$ws.Range("N107").Value = -7330.2858

# Row 108 (Leve Item ID 27853)
$ws.Range("H108").Value = 1979.2
$ws.Range("I108").Value = 474
$ws.Range("J108").Value = 8000
$ws.Range("K108").Value = 1422
$ws.Range("L108").Value = 24000
$ws.Range("M108").Value = 1458
$ws.Range("N108").Value = -29760

# Row 109 (Leve Item ID 27854)
$ws.Range("H109").Value = 818.8333
$ws.Range("I109").Value = 463.5
$ws.Range("J109").Value = 996.5
$ws.Range("K109").Value = 1390.5
$ws.Range("L109").Value = 2989.5
$ws.Range("M109").Value = -350.5
$ws.Range("N109").Value = -5069.5

# Row 114 (Leve Item ID 27865)
$ws.Range("H114").Value = 1098
$ws.Range("I114").Value = 1123.75
$ws.Range("J114").Value = 995
$ws.Range("K114").Value = 3371.25
$ws.Range("L114").Value = 2985
$ws.Range("M114").Value = -117.25
$ws.Range("N114").Value = -9493

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 1324
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 1432
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 4296
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -14376


# ============ Sheet: GSM ============
$ws = $wb.Worksheets.Item("GSM")

# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 2503
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 2503
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

# Row 95 (Leve Item ID 18235)
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 2100
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2100
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6440

# Row 133 (Leve Item ID 41854)
$ws.Range("H133").Value = 45000
$ws.Range("J133").Value = 45000
$ws.Range("L133").Value = 45000
$ws.Range("N133").Value = -55120


# ============ Sheet: LTW ============
$ws = $wb.Worksheets.Item("LTW")

# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 1255
$ws.Range("I16").Value = 1399
$ws.Range("J16").Value = 1111
$ws.Range("K16").Value = 1399
$ws.Range("L16").Value = 1111
$ws.Range("M16").Value = -1229
$ws.Range("N16").Value = -1451

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 499.5
$ws.Range("I22").Value = 499.5
$ws.Range("K22").Value = 499.5
$ws.Range("M22").Value = -204.5

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 499.5
$ws.Range("I27").Value = 499.5
$ws.Range("K27").Value = 499.5
$ws.Range("M27").Value = -392.5

# Row 43 (Leve Item ID 4314)
$ws.Range("H43").Value = 15000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

# Row 45 (Leve Item ID 3851)
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 2745.3572
$ws.Range("I100").Value = 1492.2858
$ws.Range("K100").Value = 1492.2858
$ws.Range("M100").Value = -951.2858000000001

# Row 104 (Leve Item ID 18675)
$ws.Range("H104").Value = 12500
$ws.Range("J104").Value = 12500
$ws.Range("L104").Value = 12500
$ws.Range("N104").Value = -19488

# Row 106 (Leve Item ID 18713)
$ws.Range("H106").Value = 25166.334
$ws.Range("J106").Value = 25166.334
$ws.Range("L106").Value = 25166.334
$ws.Range("N106").Value = -27690.334

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 2516.7144
$ws.Range("I132").Value = 2516.7144
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7550.1432
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5020.1432
$ws.Range("N132").ClearContents()

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 3188
$ws.Range("I136").Value = 2749.8333
$ws.Range("K136").Value = 8249.499899999999
$ws.Range("M136").Value = -5699.499899999999


# ============ Sheet: WVR ============
$ws = $wb.Worksheets.Item("WVR")

# Row 29 (Leve Item ID 3568)
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("K29").Value = 1000
$ws.Range("M29").Value = -710

# Row 82 (Leve Item ID 10806)
$ws.Range("H82").Value = 34750
$ws.Range("J82").Value = 34750
$ws.Range("L82").Value = 34750
$ws.Range("N82").Value = -35516

# Row 85 (Leve Item ID 10806)
$ws.Range("H85").Value = 34750
$ws.Range("J85").Value = 34750
$ws.Range("L85").Value = 34750
$ws.Range("N85").Value = -37402

# Row 96 (Leve Item ID 19977)
$ws.Range("H96").Value = 2062
$ws.Range("I96").Value = 1480.625
$ws.Range("K96").Value = 1480.625
$ws.Range("M96").Value = -107.625

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 38309.15
$ws.Range("I132").Value = 42889.875
$ws.Range("J132").Value = 1663.3334
$ws.Range("K132").Value = 128669.625
$ws.Range("L132").Value = 4990.0002
$ws.Range("M132").Value = -126139.625
$ws.Range("N132").Value = -10050.0002
